# Weekly fruit/vegetable price update.
# A new weekly record is inserted at row 137 (pushing the existing rows
# 137-149 down to 138-150), and the new row is populated with this week's
# figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 137; this shifts rows 137:149 down to 138:150
# and carries the row's number formatting (date style on column D) down
# with it, matching native Excel "Insert Row" behaviour.
$ws.Rows("137").Insert()

# Populate the newly inserted row 137 with this week's data. Most fields
# repeat the constant attributes of this market/product subset; only the
# date, volume, prices and per-kg price are new for this observation.
$ws.Cells.Item(137, 1).Value = 10
$ws.Cells.Item(137, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(137, 3).Value = "La Araucanía"
$ws.Cells.Item(137, 4).Value = 44449
$ws.Cells.Item(137, 5).Value = 9
$ws.Cells.Item(137, 6).Value = 100112043
$ws.Cells.Item(137, 7).Value = "Pepino dulce"
$ws.Cells.Item(137, 8).Value = "Cultivar IV Región"
$ws.Cells.Item(137, 9).Value = "Primera"
$ws.Cells.Item(137, 10).Value = 45
$ws.Cells.Item(137, 11).Value = 22000
$ws.Cells.Item(137, 12).Value = 22000
$ws.Cells.Item(137, 13).Value = 22000
$ws.Cells.Item(137, 14).Value = "$/bandeja 18 kilos"
$ws.Cells.Item(137, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(137, 16).Value = 1222
$ws.Cells.Item(137, 17).Value = 18
$ws.Cells.Item(137, 18).Value = "Hortaliza"
